$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct a tiny floating-point precision difference in the existing last row (A32)
$ws.Cells.Item(32, 1).Value = 44345.80937171412

# Append the newly retrieved row of data
$ws.Cells.Item(33, 1).Value = 44346.80724647864
$ws.Cells.Item(33, 2).Value = 73794
$ws.Cells.Item(33, 3).Value = 62079
$ws.Cells.Item(33, 4).Value = 3313
$ws.Cells.Item(33, 5).Value = 2036
$ws.Cells.Item(33, 6).Value = 1423
$ws.Cells.Item(33, 7).Value = 19209
$ws.Cells.Item(33, 8).Value = 1375
$ws.Cells.Item(33, 9).Value = 813
$ws.Cells.Item(33, 10).Value = 199

# Match the date style used in column A (style index 2 -> numFmtId 164)
$ws.Range("A33").NumberFormat = $ws.Range("A32").NumberFormat
